$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, pushing existing rows 23..40 down to 24..41
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly record
$ws.Cells.Item(23, 1).Value = 1
$ws.Cells.Item(23, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(23, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(23, 4).Value = 44566
$ws.Cells.Item(23, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23, 5).Value = 15
$ws.Cells.Item(23, 6).Value = 100112027
$ws.Cells.Item(23, 7).Value = "Melón"
$ws.Cells.Item(23, 8).Value = "Calameño"
$ws.Cells.Item(23, 9).Value = "Segunda"
$ws.Cells.Item(23, 10).Value = 50
$ws.Cells.Item(23, 11).Value = 6000
$ws.Cells.Item(23, 12).Value = 7000
$ws.Cells.Item(23, 13).Value = 6500
$ws.Cells.Item(23, 14).Value = "$/caja 24 unidades"
$ws.Cells.Item(23, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(23, 16).Value = 271
$ws.Cells.Item(23, 17).Value = 24
$ws.Cells.Item(23, 18).Value = "Hortaliza"
